$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new BOM line (for the Raspberry Pi itself) was inserted right after row 2,
# pushing the rest of the (empty, formula-only) template rows down by one.
$ws.Range("A3").EntireRow.Insert()

# --- New BOM rows (2-7) content -------------------------------------------------
$ws.Range("B2").Value = "Raspberry Pi 3 Model B Vi. 2 "
$ws.Range("C2").Value = "Raspberry Pi "
$ws.Range("D2").Value = "RASPBERRY PI 3"
$ws.Range("E2").Value = "-"
$ws.Range("F2").Value = "-"
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = "https://www.raspberrypi.com/products/raspberry-pi-3-model-b/"

$ws.Range("B3").Value = "Micro SD Card - 16 GB - with adapter"
$ws.Range("C3").Value = "Sandisk"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "Adafruit"
$ws.Range("F3").Value = 2693
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 19.95
$ws.Range("I3").Formula = "=G3*H3"
$ws.Range("J3").Value = "https://www.adafruit.com/product/2693"

$ws.Range("B4").Value = "Raspberry Pi Camera Board v2 - 8 Megapixels"
$ws.Range("C4").Value = "Raspberry Pi "
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "Adafruit"
$ws.Range("F4").Value = 3099
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 29.95
$ws.Range("J4").Value = "https://www.adafruit.com/product/3099"

$ws.Range("B5").Value = "Micro B USB Cable w/ LCD Voltage/Current Display"
$ws.Range("C5").Value = "Adafruit "
$ws.Range("D5").Value = 3388
$ws.Range("E5").Value = "Adafruit "
$ws.Range("F5").Value = 3388
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 9.95
$ws.Range("J5").Value = "https://www.adafruit.com/product/3388"

$ws.Range("B6").Value = "USB WiFi Module"
$ws.Range("C6").Value = "OURLink"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "Adafruit"
$ws.Range("F6").Value = 1012
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 13.95
$ws.Range("J6").Value = "https://www.adafruit.com/product/1012"

$ws.Range("B7").Value = "PIR Motion Sensor "
$ws.Range("C7").Value = "Adafruit "
$ws.Range("D7").Value = 189
$ws.Range("E7").Value = "Adafruit"
$ws.Range("F7").Value = 189
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 9.95
$ws.Range("J7").Value = "https://www.adafruit.com/product/189"
